$d = $word.ActiveDocument

# 1) Nomenclature label updates (simple text replacements)
$d.Content.Find.Execute("(NOMENCLATURE: ENUM-ROLE)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(NOMENCLATURE: HubSante.role)", 2) | Out-Null

$d.Content.Find.Execute("(NOMENCLATURE: NOS-NOMENC_SEXE)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(NOMENCLATURE: HubSante.sexe)", 2) | Out-Null

$d.Content.Find.Execute("(NOMENCLATURE: ENUM-TYPE_Id_Patient)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(NOMENCLATURE: HubSante.typeIdPatient)", 2) | Out-Null

# 2) Add a new "precision" row right after the "value" row in the last table
$t = $d.Tables.Item(6)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "precision"
$newRow.Cells.Item(2).Range.Text = "Précision sur la mesure"
$newRow.Cells.Item(3).Range.Text = "string"
$newRow.Cells.Item(4).Range.Text = "0..1"
$newRow.Cells.Item(5).Range.Text = "bras droit/gauche, débit oxygène, …"
